$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D used to hold boolean FALSE values. It now holds text
# representations of a date (with distinguishing trailing spaces so the two
# rows end up as distinct shared-string entries), formatted the same way as
# the date columns B/C. Copy C1's format (via PasteSpecial so the existing
# style index is reused instead of a new numFmt being minted) onto D1:D2,
# then set the text values.
$ws.Range("C1").Copy()
$ws.Range("D1:D2").PasteSpecial(-4122)
$ws.Range("D1").Value = "14/12/2021  "
$ws.Range("D2").Value = "14/12/2021 "

# Column D should now be about as wide as column C.
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(3).ColumnWidth

# New column E holds the boolean values that used to live in D.
$ws.Range("E1").Value = $false
$ws.Range("E2").Value = $false

# Update the last-used selection to reflect the new layout.
$ws.Range("D4").Select() | Out-Null
